# Auto-generated edit script: update D (Price) and E (Volume(1h)) columns
# to match the refreshed cryptos snapshot, preserving each cell's text type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.957.55'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.818.33'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4692'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3668'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07352'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8731'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '1.816.28'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.411'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07114'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.513'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008709'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = '26.971.30'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.293'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("D24").Value = '2.041.39'
$ws.Range("E24").Value = '  -0.91%  '
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.154'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("E29").Value = '  -0.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08901'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7613'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.14%  '
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.505'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.909'
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.096'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05295'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.966'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.395'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5294'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.156'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1655'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.449'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4875'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.81%  '
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06295'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.26%  '
